$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at H (shifts old H "item_specialty_reward_type" and
#    I "artifact_item_id" one column to the right, to I and J respectively).
$ws.Columns.Item(8).Insert()

# 2. Give the new column H a header and set its width to (approximately) the
#    target bestFit width for the long descriptive text it will hold.
$ws.Range("H1").Value2 = "scheduled_event_description"
$ws.Range("H1").EntireColumn.ColumnWidth = 527.61

# 3. Populate the new scheduled_event_description cells for the two raids
#    that received flavor text for their in-game scheduled events.
$ws.Range("H3").Value2 = "Take part in a raid against the queen her self. She looks for her son, she grieves her husband. She is the queen of this realm and her suffering needs to end. Participate in The Ice Queens Reign Raid to earn a cosmetic item through the raid quests, and ancestral item for being the first to defeat her and a full set of Corrupted Ice gear (a more powerful set of gear!) for defeating her! All players are welcome to try their strength against her!"
$ws.Range("H5").Value2 = "The king seeks his son, but he cries out for his wife. His is created of corrupted tears and failed magics. Take him down to earn an ancestral item (as the first person to kill him) and a full set of Corrupted Ice gear (a more powerful set of gear) as well! try your strength against him, regardless of your level!."

# 4. The raid_monster_ids list for The Frozen King raid (row 5) was reordered.
$ws.Range("E5").Value2 = "Faithless Priest of The Old Church,Zombified Cat of Yesterday,Faithless Prince of the Snow Garden,Bloody Snowman of rage,Wailing Banshee of Ice,Corrupted Christmas Tree,Frozen Child of Fear"
